$d = $word.ActiveDocument

# Smart quote / apostrophe characters used throughout the document.
$rsquo = [char]8217
$ldquo = [char]8220
$rdquo = [char]8221

# ------------------------------------------------------------------
# Helper: wipe out *all* the content of a paragraph (including any
# extra runs / proofErr markers it may contain) while leaving the
# paragraph mark itself (and therefore the paragraph + its pPr/style)
# intact, so a fresh run can be inserted into it afterwards.
# ------------------------------------------------------------------
function Clear-ParaText([int]$index) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    if (($r.End - 1) -gt $r.Start) {
        $sub = $d.Range($r.Start, $r.End - 1)
        $sub.Delete()
    }
}

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sat in the empty Image
#    paragraph near the top of the document.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Rewrite the "Challenge 1" / "Challenge 2" block.
#    Paragraphs 31-38 (pre-edit) become 5 paragraphs:
#      31 Heading2  -> "Challenge 1: " + "Bring Back the Aliens!"
#      32 Normal    -> alien update text
#      33 Heading2  -> "Challenge 2: " + "Refactor, Refactor, Refactor!"
#      34 Normal    -> refactor text
#      35 Normal    -> sections/rows text (+ _GoBack bookmark at end)
#    Paragraphs 36-38 (old leftover content) are removed.
# ------------------------------------------------------------------

# --- Paragraph 31: Heading2 "Challenge 1: Editing Your Stories" ---
Clear-ParaText 31
$d.Paragraphs.Item(31).Range.InsertBefore("Challenge 1: ")
$d.Paragraphs.Item(31).Range.InsertAfter("Bring Back the Aliens!")

# --- Paragraph 32: body text under Challenge 1 ---
Clear-ParaText 32
$alienText = "If you" + $rsquo + "ve been following along since the first session, you" + $rsquo + "ll remember that there were Aliens and now there is not. Update the app so that there is now alien stories to choose from. Make sure to include your own alien story."
$d.Paragraphs.Item(32).Range.InsertBefore($alienText)

# --- Paragraph 33: was Normal ("Your challenge is..."), becomes Heading2 "Challenge 2: Refactor, Refactor, Refactor!" ---
Clear-ParaText 33
$d.Paragraphs.Item(33).Range.InsertBefore("Challenge 2: ")
$d.Paragraphs.Item(33).Range.InsertAfter("Refactor, Refactor, Refactor!")
$d.Paragraphs.Item(33).Style = "Heading 2"

# --- Paragraph 34: was Normal ("To let the user know..."), becomes refactor body text ---
Clear-ParaText 34
$refactorText = "Currently, the app is using individual arrays for each type of story. Refactor the app to use one array. This array will contain three arrays: one for the zombies, one for the vampires, and one for the aliens."
$d.Paragraphs.Item(34).Range.InsertBefore($refactorText)

# --- Paragraph 35: was Normal ("Once the user taps this accessory..."), becomes the sections/rows sentence, with the _GoBack bookmark placed at its end ---
Clear-ParaText 35
$sectionsText = "You should be able to access the individual stories by using sections and rows"
$d.Paragraphs.Item(35).Range.InsertBefore($sectionsText)
$p35end = $d.Paragraphs.Item(35).Range.End
$d.Bookmarks.Add("_GoBack", $d.Range($p35end, $p35end)) | Out-Null

# --- Paragraphs 36, 37, 38 (old "That means...", "Challenge 2: Deleting Your Stories", "Let's face it...") are no longer needed ---
$d.Paragraphs.Item(38).Range.Delete()
$d.Paragraphs.Item(37).Range.Delete()
$d.Paragraphs.Item(36).Range.Delete()
